$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Encoded per-row category codes for rows 2..451 (one char per row):
#   O = "Operating Activities (XGBoost)"
#   I = "Investing Activities (XGBoost)"
#   F = "Financing Activities (XGBoost)"
$codes = "IFIFIOOOOIIIOOFOOOOIIOFOIFOIOOOIIIIIOIIIOOIOIIOIOIFFOIFFIIIOIOFIIFOIOIIOIOOIIOOIIIOIIIOOOIIIOFOIOIFOOOOIOIOOIIOIOIIIIOOOFOOIIOFIIOOIIFOOOOOOIOIIIIOOOIIOOIOFIOOOOIOIOIOOIIIOOIIOIOIIOOOIIFIIFOIFFOFIIIOIOOFIIIOIOIOIFFIFOIFOIIIOOOIOOIIIOIOOOFOIIOIOIIOOOIOIIOOOIOOIIOOFFOIIOIFIIOOIIOOOIOIIIOOFOFIOIOOIIIIOOIIOOOIIFOIOOIFIOIIOOFIOOIOFFOOOFOOIFIOOIOIOOOIOIOOOOIIOOOIOOIOOIOIOFOOFOFIOFOOIOOIIIOOOIIIOIOIIOIIIFIIIOIIIIIIOOOOIIIOFIIOFIFOOIOOIIIOIIOOOOOIIFIOFFO"

$labels = @{
    "O" = "Operating Activities (XGBoost)"
    "I" = "Investing Activities (XGBoost)"
    "F" = "Financing Activities (XGBoost)"
}

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $code = $codes.Substring($i, 1)
    $label = $labels[$code]
    $ws.Range("F" + $row).Value = $label
}

Write-Output "Updated $($codes.Length) rows in column F"
